$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.262227177619934
$ws.Range("B1").Value = 2.494179964065552
$ws.Range("C1").Value = 4.995416164398193
$ws.Range("D1").Value = 2.926954984664917
$ws.Range("E1").Value = 1.105132102966309
